# "Made arcing shots more natural looking"
# Adjust the base pass-parabola control points (B2 start height, D2 apex
# height) and add a small net-clearance calculator in columns H:I that
# derives how much the apex height needs to be nudged so the arc still
# clears the net by the required margin.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Core parabola control points -----------------------------------
$ws.Range("B2").Value = 3.5
$ws.Range("D2").Value = 3.62

# --- New net-clearance helper block (columns H/I) --------------------
$ws.Range("H6").Value = "Start-x"
$ws.Range("I6").Value = -3.5

$ws.Range("H7").Value = "End-x"
$ws.Range("I7").Value = 6.78

$ws.Range("H8").Value = "netCrossingT"
$ws.Range("I8").Formula = "=ABS(I6)/ABS(I7-I6)"

$ws.Range("H10").Value = "heightAtNet"
$ws.Range("I10").Formula = "=(1-I8)*(1-I8)*B2+2*(1-I8)*I8*((D2-0.25*B2-0.25*C2)/0.5)+I8*I8*C2"

$ws.Range("H11").Value = "maxHeightPoint"
$ws.Range("I11").Formula = "=D2"

$ws.Range("H12").Value = "requiredNetHeight"
$ws.Range("I12").Value = 2.5

$ws.Range("H13").Value = "adjustedMaxHeight"
$ws.Range("I13").Formula = "=I11+I12-I10"

# --- Formatting --------------------------------------------------------
# Editable "input" cells get a bold blue font so they stand out.
$inputCells = @("B2", "D2", "I6", "I7", "I12")
foreach ($addr in $inputCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Font.Color = 15773696
}

# netCrossingT mirrors the workbook's existing "0.000" display format.
$ws.Range("I8").NumberFormat = $ws.Range("E2").NumberFormat

# Key derived outputs (heightAtNet, adjustedMaxHeight) get a 2-decimal
# display format - applied after the font tweaks above so the style
# table lays out the same way it does when a person formats the inputs
# first and the summary numbers second.
$ws.Range("I10").NumberFormat = "0.00"
$ws.Range("I13").NumberFormat = "0.00"

# Column widths: a narrow spacer column (G) and a wide label column (H).
$ws.Columns.Item(7).ColumnWidth = 1.69140625
$ws.Columns.Item(8).ColumnWidth = 16.53515625

$wb.Application.Calculate()
